$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Summary sheet - update aggregate metrics after trade #73 closes
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1200.12   # Current Capital
$summary.Range("B4").Value = 0.11      # Total P&L $
$summary.Range("B5").Value = 0.03      # Total P&L %
$summary.Range("B6").Value = 73        # Total Trades
$summary.Range("B8").Value = 39        # Losing Trades
$summary.Range("B9").Value = 31.51     # Win Rate %

# ---------------------------------------------------------------------------
# 2) Strategy Status sheet - update MarketMaking row (row 4)
# ---------------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 100.12     # Capital
$status.Range("D4").Value = 73         # Trades
$status.Range("E4").Value = 0.11       # P&L $
$status.Range("F4").Value = 0.12       # P&L %
$status.Range("G4").Value = 31.51      # Win Rate %

# ---------------------------------------------------------------------------
# 3) Append the newly closed trade (#73) to both "All Trades" and
#    "MarketMaking" sheets as row 74.
# ---------------------------------------------------------------------------
$sheetsToAppend = @("All Trades", "MarketMaking")

foreach ($sheetName in $sheetsToAppend) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Cells.Item(74, 1).Value = 73                 # Trade #

    # Date / Time columns look like dates/times, but must stay plain text
    # (matching the rest of the column), so force a text number format
    # before assigning them to avoid automatic date/time conversion.
    $ws.Cells.Item(74, 2).NumberFormat = "@"
    $ws.Cells.Item(74, 2).Value = "2026-02-17"        # Date
    $ws.Cells.Item(74, 3).NumberFormat = "@"
    $ws.Cells.Item(74, 3).Value = "15:48:40"          # Time

    $ws.Cells.Item(74, 4).Value = "MarketMaking"      # Strategy
    $ws.Cells.Item(74, 5).Value = "UP"                # Side
    $ws.Cells.Item(74, 6).Value = 0.79                # Entry Price
    $ws.Cells.Item(74, 7).Value = 0.66                # Exit Price
    $ws.Cells.Item(74, 8).Value = "CLOSED"            # Status
    $ws.Cells.Item(74, 9).Value = -16.4557            # P&L %
    $ws.Cells.Item(74, 10).Value = -0.13              # P&L $
    $ws.Cells.Item(74, 11).Value = 100.12             # Capital After
    $ws.Cells.Item(74, 12).Value = 0                  # Entry Slippage (bps)
    $ws.Cells.Item(74, 13).Value = 0                  # Exit Slippage (bps)
    $ws.Cells.Item(74, 14).Value = 0.6                # Confidence
    $ws.Cells.Item(74, 15).Value = "Normal spread capture: 19600 bps"  # Entry Reason
    $ws.Cells.Item(74, 16).Value = "early_exit"       # Exit Reason
    $ws.Cells.Item(74, 17).Value = 0.15               # Duration (min)
}
